$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 3) to the meta-sheet for PF/1.0.6
$row = $ws.Range("A3:D3")
$row.Style = "Normal"

$ws.Range("A3").Value = "PF/1.0.6"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
